$d = $word.ActiveDocument

# 1. Bump the GroupDocs.Assembly version stamp printed in the evaluation
#    watermark paragraph: 25.6. -> 25.12. (December Christmas release).
$d.Content.Find.Execute(
    "Created with GroupDocs.Assembly 25.6. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Created with GroupDocs.Assembly 25.12. ", 2)

# 2. Register the (previously absent) built-in "Hyperlink" character style so
#    that it is present in styles.xml, based on Default Paragraph Font, with
#    the usual blue/underlined look (RGB 0563C1, single underline) and
#    ui priority 99, matching Word's standard Hyperlink style definition.
$hyperlinkStyle = $d.Styles.Add("Hyperlink", 2)
$hyperlinkStyle.BaseStyle = "DefaultParagraphFont"
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.Font.Color = 12673797
$hyperlinkStyle.Font.Underline = 1
